# Update the raw "x (in)" readings in both data sets to the corrected
# position values (in mm) recorded for this run, and fix a data-entry
# typo in the first set's pressure reading (a missing decimal point).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set 1: Pressures Along the Tube (rows 8-19, column E) ---
$set1 = @{
    8  = 12.75
    9  = 28.25
    10 = 52.25
    11 = 79.5
    12 = 103.5
    13 = 127.25
    14 = 151.25
    15 = 175.25
    16 = 212.25
    17 = 219.25
    18 = 229.25
    19 = 251.25
}
foreach ($row in $set1.Keys) {
    $ws.Range("E$row").Value = $set1[$row]
}

# Fix typo: 162 -> 1.62 (decimal point was dropped when recording)
$ws.Range("F15").Value = 1.62

# --- Set 2: Pressures Along the Tube (rows 24-35, column E) ---
$set2 = @{
    24 = 12.75
    25 = 28.25
    26 = 52.25
    27 = 79.5
    28 = 103.5
    29 = 127.25
    30 = 151.25
    31 = 175.25
    32 = 212.25
    33 = 219.25
    34 = 229.25
    35 = 251.25
}
foreach ($row in $set2.Keys) {
    $ws.Range("E$row").Value = $set2[$row]
}

# Update the active cell selection to reflect where the user left off.
$ws.Range("E2").Select() | Out-Null
